$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 15
    $ws.Range("F4").Value = 1463
    $ws.Range("F9").Value = 251
}
